$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F ("想去人数" / number interested) updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 525
$wsExpo.Range("F3").Value = 6303
$wsExpo.Range("F4").Value = 395
$wsExpo.Range("F5").Value = 94
$wsExpo.Range("F6").Value = 127
$wsExpo.Range("F8").Value = 71
$wsExpo.Range("F9").Value = 568
$wsExpo.Range("F10").Value = 43

# Sheet "全部类型" (All types) - same events, shifted rows due to extra entries
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 525
$wsAll.Range("F3").Value = 6303
$wsAll.Range("F4").Value = 395
$wsAll.Range("F6").Value = 94
$wsAll.Range("F7").Value = 127
$wsAll.Range("F10").Value = 71
$wsAll.Range("F11").Value = 568
$wsAll.Range("F12").Value = 43
